$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2626.925
$ws.Range("I98").Value = 2768.9722
$ws.Range("J98").Value = 1348.5
$ws.Range("K98").Value = 2768.9722
$ws.Range("L98").Value = 1348.5
$ws.Range("M98").Value = -1270.9722
$ws.Range("N98").Value = -4344.5

$ws.Range("H105").Value = 20000
$ws.Range("J105").Value = 20000
$ws.Range("L105").Value = 20000
$ws.Range("N105").Value = -26988

$ws.Range("H122").Value = 2626.925
$ws.Range("I122").Value = 2768.9722
$ws.Range("J122").Value = 1348.5
$ws.Range("K122").Value = 8306.9166
$ws.Range("L122").Value = 4045.5
$ws.Range("M122").Value = -5856.9166
$ws.Range("N122").Value = -8945.5

$ws.Range("H132").Value = 3925151.2
$ws.Range("I132").Value = 4833928
$ws.Range("J132").Value = 6051
$ws.Range("K132").Value = 14501784
$ws.Range("L132").Value = 18153
$ws.Range("M132").Value = -14499254
$ws.Range("N132").Value = -23213

$ws.Range("H135").Value = 47619732
$ws.Range("I135").Value = 352.77777
$ws.Range("J135").Value = 333336000
$ws.Range("K135").Value = 3174.99993
$ws.Range("L135").Value = 3000024000
$ws.Range("M135").Value = -639.9999299999999
$ws.Range("N135").Value = -3000029070

$ws.Range("H137").Value = 1027.8961
$ws.Range("I137").Value = 715.04346
$ws.Range("J137").Value = 1492.129
$ws.Range("K137").Value = 2145.13038
$ws.Range("L137").Value = 4476.387
$ws.Range("M137").Value = 404.8696199999999
$ws.Range("N137").Value = -9576.386999999999

$ws.Range("H138").Value = 1175.1718
$ws.Range("I138").Value = 542
$ws.Range("J138").Value = 1661.3572
$ws.Range("K138").Value = 1626
$ws.Range("L138").Value = 4984.071599999999
$ws.Range("M138").Value = 3514
$ws.Range("N138").Value = -15264.0716

$ws.Range("H141").Value = 477.375
$ws.Range("I141").Value = 477.375
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1432.125
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3747.875
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5982.162
$ws.Range("I32").Value = 4995.424
$ws.Range("J32").Value = 9863.333000000001
$ws.Range("K32").Value = 4995.424
$ws.Range("L32").Value = 9863.333000000001
$ws.Range("M32").Value = -4708.424
$ws.Range("N32").Value = -10437.333

$ws.Range("H45").Value = 1214.2858
$ws.Range("I45").Value = 1333.3334
$ws.Range("J45").Value = 1125
$ws.Range("K45").Value = 1333.3334
$ws.Range("L45").Value = 1125
$ws.Range("M45").Value = -956.3334
$ws.Range("N45").Value = -1879

$ws.Range("H61").Value = 52632628
$ws.Range("I61").Value = 62500948
$ws.Range("K61").Value = 62500948
$ws.Range("M61").Value = -62500736

$ws.Range("H74").Value = 862.4666999999999
$ws.Range("I74").Value = 683.61536
$ws.Range("K74").Value = 683.61536
$ws.Range("M74").Value = 190.38464

$ws.Range("H77").Value = 862.4666999999999
$ws.Range("I77").Value = 683.61536
$ws.Range("K77").Value = 3418.0768
$ws.Range("M77").Value = 949.9232000000002

$ws.Range("H132").Value = 1565.6078
$ws.Range("I132").Value = 1256.1945
$ws.Range("J132").Value = 2308.2
$ws.Range("K132").Value = 3768.5835
$ws.Range("L132").Value = 6924.599999999999
$ws.Range("M132").Value = -1238.5835
$ws.Range("N132").Value = -11984.6

$ws.Range("H136").Value = 52632628
$ws.Range("I136").Value = 62500948
$ws.Range("K136").Value = 187502844
$ws.Range("M136").Value = -187500294

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2425.4429
$ws.Range("I134").Value = 884.3889
$ws.Range("J134").Value = 7626.5
$ws.Range("K134").Value = 2653.1667
$ws.Range("L134").Value = 22879.5
$ws.Range("M134").Value = -118.1667000000002
$ws.Range("N134").Value = -27949.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2155.8333
$ws.Range("I31").Value = 2247.1
$ws.Range("J31").Value = 1699.5
$ws.Range("K31").Value = 2247.1
$ws.Range("L31").Value = 1699.5
$ws.Range("M31").Value = -1952.1
$ws.Range("N31").Value = -2289.5

$ws.Range("H34").Value = 2155.8333
$ws.Range("I34").Value = 2247.1
$ws.Range("J34").Value = 1699.5
$ws.Range("K34").Value = 2247.1
$ws.Range("L34").Value = 1699.5
$ws.Range("M34").Value = -2045.1
$ws.Range("N34").Value = -2103.5

$ws.Range("H122").Value = 1329.875
$ws.Range("J122").Value = 1195
$ws.Range("L122").Value = 3585
$ws.Range("N122").Value = -8485

$ws.Range("H132").Value = 3424.3462
$ws.Range("I132").Value = 3877.6316
$ws.Range("K132").Value = 11632.8948
$ws.Range("M132").Value = -9102.8948

$ws.Range("H134").Value = 11906181
$ws.Range("I134").Value = 1486.3103
$ws.Range("J134").Value = 38462810
$ws.Range("K134").Value = 4458.9309
$ws.Range("L134").Value = 115388430
$ws.Range("M134").Value = -1923.9309
$ws.Range("N134").Value = -115393500

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 5984
$ws.Range("I11").Value = 4313.3335
$ws.Range("J11").Value = 8490
$ws.Range("K11").Value = 12940.0005
$ws.Range("L11").Value = 25470
$ws.Range("M11").Value = -12800.0005
$ws.Range("N11").Value = -25750

$ws.Range("H81").Value = 2669.6365
$ws.Range("J81").Value = 3074.2
$ws.Range("L81").Value = 9222.599999999999
$ws.Range("N81").Value = -11468.6

$ws.Range("H84").Value = 2669.6365
$ws.Range("J84").Value = 3074.2
$ws.Range("L84").Value = 27667.8
$ws.Range("N84").Value = -38899.8

$ws.Range("H113").Value = 667.2857
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 734.2
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2202.6
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -6542.6

$ws.Range("H131").Value = 26319406
$ws.Range("J131").Value = 5053.769
$ws.Range("L131").Value = 15161.307
$ws.Range("N131").Value = -25241.307

$ws.Range("H140").Value = 21420.908
$ws.Range("I140").Value = 61809.883
$ws.Range("K140").Value = 185429.649
$ws.Range("M140").Value = -180249.649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 633.93335
$ws.Range("I102").Value = 667.4167
$ws.Range("K102").Value = 667.4167
$ws.Range("M102").Value = 954.5833

$ws.Range("H126").Value = 2129.1667
$ws.Range("I126").Value = 1718.75
$ws.Range("J126").Value = 2950
$ws.Range("K126").Value = 5156.25
$ws.Range("L126").Value = 8850
$ws.Range("M126").Value = -2686.25
$ws.Range("N126").Value = -13790

$ws.Range("H132").Value = 2289.8206
$ws.Range("J132").Value = 2836.5
$ws.Range("L132").Value = 8509.5
$ws.Range("N132").Value = -13569.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 11340776
$ws.Range("I122").Value = 16674241
$ws.Range("J122").Value = 7163.125
$ws.Range("K122").Value = 50022723
$ws.Range("L122").Value = 21489.375
$ws.Range("M122").Value = -50020273
$ws.Range("N122").Value = -26389.375

$ws.Range("H132").Value = 32559.516
$ws.Range("I132").Value = 1809.5294
$ws.Range("J132").Value = 65231.375
$ws.Range("K132").Value = 5428.5882
$ws.Range("L132").Value = 195694.125
$ws.Range("M132").Value = -2898.5882
$ws.Range("N132").Value = -200754.125

$ws.Range("H136").Value = 3074.1042
$ws.Range("I136").Value = 3206.1904
$ws.Range("J136").Value = 2149.5
$ws.Range("K136").Value = 9618.5712
$ws.Range("L136").Value = 6448.5
$ws.Range("M136").Value = -7068.5712
$ws.Range("N136").Value = -11548.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 392.17648
$ws.Range("I107").Value = 360.58334
$ws.Range("K107").Value = 1081.75002
$ws.Range("M107").Value = 838.2499800000001

$ws.Range("H113").Value = 396.57895
$ws.Range("I113").Value = 305.77777
$ws.Range("J113").Value = 478.3
$ws.Range("K113").Value = 917.33331
$ws.Range("L113").Value = 1434.9
$ws.Range("M113").Value = 1252.66669
$ws.Range("N113").Value = -5774.9

$ws.Range("H126").Value = 50506460
$ws.Range("I126").Value = 52911384
$ws.Range("K126").Value = 158734152
$ws.Range("M126").Value = -158731682

$ws.Range("H132").Value = 2213.4426
$ws.Range("I132").Value = 2109.3333
$ws.Range("J132").Value = 2443.5789
$ws.Range("K132").Value = 6327.999899999999
$ws.Range("L132").Value = 7330.736699999999
$ws.Range("M132").Value = -3797.999899999999
$ws.Range("N132").Value = -12390.7367

$ws.Range("H136").Value = 620.4103
$ws.Range("I136").Value = 544.0333000000001
$ws.Range("J136").Value = 875
$ws.Range("K136").Value = 1632.0999
$ws.Range("L136").Value = 2625
$ws.Range("M136").Value = 917.9000999999998
$ws.Range("N136").Value = -7725
